$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$ws.Range("A11").Value = "{'Hobby': '1', 'SexualOrientation': '1'}"
$ws.Range("A12").Value = "{'Hobby': '1', 'Student': '1'}"
$ws.Range("A13").Value = "{'Hobby': '1', 'Student': '1', 'SexualOrientation': '1'}"
$ws.Range("A14").Value = "{'HDI': '1', 'SexualOrientation': '1'}"
$ws.Range("A15").Value = "{'HDI': '1', 'Student': '1'}"
$ws.Range("A17").Value = "{'HDI': '1', 'Student': '1', 'SexualOrientation': '1'}"
$ws.Range("A18").Value = "{'Hobby': '1', 'SexualOrientation': '1', 'HDI': '1'}"
$ws.Range("A24").Value = "{'RaceEthnicity': '1', 'HDI': '1', 'SexualOrientation': '1'}"
$ws.Range("A25").Value = "{'RaceEthnicity': '1', 'HDI': '1', 'Student': '1'}"
$ws.Range("A26").Value = "{'RaceEthnicity': '1', 'Hobby': '1', 'SexualOrientation': '1'}"
$ws.Range("A27").Value = "{'Dependents': '2', 'SexualOrientation': '1'}"
$ws.Range("A28").Value = "{'Dependents': '2', 'Hobby': '1'}"
$ws.Range("A29").Value = "{'UndergradMajor': '2', 'SexualOrientation': '1'}"
